$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Note: assigning a numeric-looking string (e.g. "3" or "1,294") straight to
# Range.Value gets auto-coerced back into a real number by Excel's smart
# parsing. To force genuine text cells (matching the target, which stores
# these as text rather than numeric <v> cells) each cell is touched as:
#   1. NumberFormat = "@" (Text) first,
#   2. assign the string value,
#   3. Style reset back to "Normal" afterwards so no stray number format
#      sticks around on the cell (keeps formatting identical to before).
# Each cell is handled inline (not via a helper function) because passing
# COM Range objects through a PowerShell function parameter here drops
# writes intermittently.
# ---------------------------------------------------------------------------

# ===========================================================================
# Sheet "Overall": A2 1294 -> "1,294" (text)
# ===========================================================================
$wsOverall = $wb.Worksheets.Item("Overall")
$cell = $wsOverall.Range("A2")
$cell.NumberFormat = "@"
$cell.Value = "1,294"
$cell.Style = "Normal"

# ===========================================================================
# Sheet "County": B2:B97 numeric -> text (same digits, no thousands comma),
# plus a new Total row 98.
# ===========================================================================
$wsCounty = $wb.Worksheets.Item("County")

$countyCounts = @{
    2=3; 3=6; 4=8; 5=3; 6=5; 7=58; 8=12; 9=11; 10=6; 11=5; 12=2; 13=2; 14=11;
    15=6; 16=7; 17=22; 18=5; 19=5; 20=4; 21=10; 22=7; 23=22; 24=5; 25=16;
    26=2; 27=5; 28=4; 29=19; 30=5; 31=64; 32=1; 33=6; 34=5; 35=9; 36=3;
    37=2; 38=5; 39=5; 40=2; 41=5; 42=9; 43=6; 44=8; 45=5; 46=5; 47=2; 48=6;
    49=6; 50=6; 51=10; 52=72; 53=10; 54=1; 55=6; 56=9; 57=80; 58=2; 59=3;
    60=5; 61=6; 62=10; 63=12; 64=11; 65=7; 66=3; 67=3; 68=5; 69=14; 70=11;
    71=1; 72=7; 73=6; 74=12; 75=2; 76=218; 77=23; 78=8; 79=2; 80=5; 81=53;
    82=7; 83=32; 84=41; 85=4; 86=3; 87=9; 88=10; 89=10; 90=14; 91=2; 92=18;
    93=9; 94=18; 95=58; 96=2; 97=4
}

foreach ($r in $countyCounts.Keys) {
    $cell = $wsCounty.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = [string]$countyCounts[$r]
    $cell.Style = "Normal"
}

# New Total row (row 98)
$cell = $wsCounty.Range("A98")
$cell.NumberFormat = "@"
$cell.Value = "Total"
$cell.Style = "Normal"

$cell = $wsCounty.Range("B98")
$cell.NumberFormat = "@"
$cell.Value = "1,294"
$cell.Style = "Normal"

$cell = $wsCounty.Range("C98")
$cell.NumberFormat = "@"
$cell.Value = "$1,252,625,535"
$cell.Style = "Normal"

$cell = $wsCounty.Range("D98")
$cell.NumberFormat = "@"
$cell.Value = "10.37%"
$cell.Style = "Normal"

$cell = $wsCounty.Range("E98")
$cell.NumberFormat = "@"
$cell.Value = "-5.73%"
$cell.Style = "Normal"

$cell = $wsCounty.Range("F98")
$cell.NumberFormat = "@"
$cell.Value = "62.36%"
$cell.Style = "Normal"

# ===========================================================================
# Sheet "Congressional District": B2:B5 numeric -> text, B6 (Total) -> "1,294"
# ===========================================================================
$wsCd = $wb.Worksheets.Item("Congressional District")

$cdCounts = @{2="291"; 3="339"; 4="319"; 5="345"; 6="1,294"}
foreach ($r in $cdCounts.Keys) {
    $cell = $wsCd.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $cdCounts[$r]
    $cell.Style = "Normal"
}

# ===========================================================================
# Sheet "Size": B2:B7 numeric -> text, B8 (Total) -> "1,294"
# ===========================================================================
$wsSize = $wb.Worksheets.Item("Size")

$sizeCounts = @{2="440"; 3="275"; 4="185"; 5="106"; 6="188"; 7="100"; 8="1,294"}
foreach ($r in $sizeCounts.Keys) {
    $cell = $wsSize.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $sizeCounts[$r]
    $cell.Style = "Normal"
}

# ===========================================================================
# Sheet "Subsector": B2:B13 numeric -> text, B14 (Total) -> "1,294"
# ===========================================================================
$wsSub = $wb.Worksheets.Item("Subsector")

$subCounts = @{2="117"; 3="71"; 4="41"; 5="118"; 6="23"; 7="472"; 8="6"; 9="1"; 10="101"; 11="25"; 12="291"; 13="28"; 14="1,294"}
foreach ($r in $subCounts.Keys) {
    $cell = $wsSub.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $subCounts[$r]
    $cell.Style = "Normal"
}
